$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename header row strings from "..._old" -> "..._FV2210" and
#    "..._new" -> "..._FV2304" (columns A-J are the "_old" block, K is
#    "diff" and stays as-is, L-U are the "_new" block).
# ---------------------------------------------------------------------
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = ($cell.Value2 -replace "_old$", "_FV2210")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = ($cell.Value2 -replace "_new$", "_FV2304")
}

# ---------------------------------------------------------------------
# 2) Freeze the header row (top row frozen, same as Excel's
#    "Freeze Top Row" command).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# 3) Turn the data range into a native Excel Table ("Table1") so the
#    header row gets filter buttons and structured references work.
#    The header row already carries explicit formatting (bold/fill/
#    border) from the sheet's own styles; ListObjects.Add would
#    otherwise bake a duplicate "header style" into styles.xml as a new
#    dxf. To avoid mutating styles.xml, stash the header formatting,
#    reset it to the default "Normal" style before creating the table,
#    then restore it once the table exists.
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A200:U200")

$headerRange.Copy()
$scratchRange.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.Style = "Normal"

$dataRange = $ws.Range("A1:U92")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)    # xlPasteFormats
$excel.CutCopyMode = $false

$scratchRange.Clear()
